# "Fixed Formatting" — see commit message.
#
# Summary of the edit being reproduced:
#  1. The active tab moves from "Seventh Astral Era MSQ" to "2.5".
#  2. View/selection state changes on several sheets (cosmetic, matches the
#     author re-saving after clicking around).
#  3. "2.5 MSQ": quest-name cells were center-aligned; they become left-aligned.
#  4. "2.55 MSQ" ("Before the Fall - Part 2 Quests") gets the same left-align
#     fix, PLUS it picks up the "Level" (column E) values that the other
#     quest sheets already had (every row = 50), and the two rows for
#     "Ancient Ways, Timeless Wants" / "Where We Are Needed" get corrected:
#     the long quest title now properly spans a merged two-row cell
#     (A12:A13 / B12:D13) instead of colliding with the next entry.

$wb = $excel.ActiveWorkbook

$xlLeft   = -4131
$xlCenter = -4108

# ---------------------------------------------------------------------
# Sheet "2.5 MSQ": quest-title cells go from centered to left-aligned.
# ---------------------------------------------------------------------
$wsMsq25 = $wb.Worksheets.Item("2.5 MSQ")
$wsMsq25.Range("B3:D10").HorizontalAlignment = $xlLeft

# ---------------------------------------------------------------------
# Sheet "2.55 MSQ": fix alignment + fill in missing Level (E) values +
# fix the "Ancient Ways, Timeless Wants" / "Where We Are Needed" rows.
# ---------------------------------------------------------------------
$wsMsq255 = $wb.Worksheets.Item("2.55 MSQ")

# Rows 3-11 keep their row numbers; just merge B:D (if not already), align
# left and stamp the Level value of 50 that every sister-sheet already has.
foreach ($r in 3..11) {
    $rng = $wsMsq255.Range("B$r`:D$r")
    if (-not $rng.MergeCells) { $rng.Merge() }
    $rng.HorizontalAlignment = $xlLeft
    $wsMsq255.Range("E$r").Value = 50
}

# Insert a row below row 12 ("Ancient Ways, Timeless Wants") so that it can
# be combined with what is currently row 13 ("Where We Are Needed") into a
# single, taller, wrapped cell — and push "Where We Are Needed" (and
# everything after it) down by one row.
$wsMsq255.Rows.Item(13).Insert()

# Row 12+13 now form the "Ancient Ways, Timeless Wants" entry.
$wsMsq255.Range("A12:A13").Merge()
$wsMsq255.Range("A12:A13").HorizontalAlignment = $xlCenter
$wsMsq255.Range("A12:A13").VerticalAlignment = $xlCenter

$wsMsq255.Range("B12:D13").Merge()
$wsMsq255.Range("B12:D13").HorizontalAlignment = $xlLeft
$wsMsq255.Range("B12:D13").WrapText = $true

$wsMsq255.Range("E12").Value = 50
$wsMsq255.Range("E13").Value = 50

# Rows that used to be 13-18 are now 14-19; they need the same treatment as
# rows 3-11 above (merge / left-align / Level = 50).
foreach ($r in 14..19) {
    $rng = $wsMsq255.Range("B$r`:D$r")
    if (-not $rng.MergeCells) { $rng.Merge() }
    $rng.HorizontalAlignment = $xlLeft
    $wsMsq255.Range("E$r").Value = 50
}

# ---------------------------------------------------------------------
# View / selection state.
# ---------------------------------------------------------------------

# "Seventh Astral Era MSQ": selection shrinks from A9:E13 down to A10.
$wsSeventh = $wb.Worksheets.Item("Seventh Astral Era MSQ")
$wsSeventh.Select()
$wsSeventh.Range("A10").Select()

# "2.5 MSQ": selection moves from E11 to B11.
$wsMsq25.Range("B11").Select()

# "2.55 MSQ": selection moves to E22 (below/right of the used range).
$wsMsq255.Range("E22").Select()

# "2.5" becomes the active tab (was "Seventh Astral Era MSQ").
$wsHome = $wb.Worksheets.Item("2.5")
$wsHome.Activate()
$wsHome.Range("A1").Select()
